$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A59").Value = "2025/12/05 04:00"
$ws.Range("B59").Value = "14,687位本"
$ws.Range("C59").Value = "45位 広告・宣伝 (本)"
$ws.Range("D59").Value = "55位商業デザイン"
$ws.Range("E59").Value = "931位ビジネス実用本"
$ws.Range("F59").Value = "-"
$ws.Range("G59").Value = "-"
